# Trade #72 closed at 2026-02-17 08:57:43 - unknown UNKNOWN +0.000%
#
# Updates the rolling trading-results workbook with the newly closed
# MarketMaking trade (#72 / 0-based row 73):
#   - Summary sheet: refreshed aggregate stats.
#   - Strategy Status sheet: refreshed MarketMaking strategy row.
#   - All Trades / MarketMaking sheets: append the new trade row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.56   # Current Capital
$summary.Range("B4").Value = 0.57      # Total P&L $
$summary.Range("B5").Value = 0.16      # Total P&L %
$summary.Range("B6").Value = 72        # Total Trades
$summary.Range("B7").Value = 30        # Winning Trades
$summary.Range("B9").Value = 41.67     # Win Rate %

# ---------------------------------------------------------------------
# 2) Strategy Status sheet - MarketMaking row (row 4)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.56     # Capital
$status.Range("D4").Value = 72         # Trades
$status.Range("E4").Value = 0.57       # P&L $
$status.Range("F4").Value = 0.56       # P&L %
$status.Range("G4").Value = 41.67      # Win Rate %

# ---------------------------------------------------------------------
# 3) Append the new closed trade (trade #72) as row 73 on both the
#    "All Trades" ledger and the per-strategy "MarketMaking" ledger.
# ---------------------------------------------------------------------
function Add-TradeRow73($ws) {
    $ws.Cells.Item(73, 1).Value = 72

    # Date / Time columns look like dates to COM's type-inference, but the
    # sheet stores them as plain text - force text format, assign, then
    # drop the format again so no stray style survives on the cell.
    $ws.Cells.Item(73, 2).NumberFormat = "@"
    $ws.Cells.Item(73, 2).Value = "2026-02-17"
    $ws.Cells.Item(73, 2).ClearFormats()

    $ws.Cells.Item(73, 3).Value = "08:57:37"

    $ws.Cells.Item(73, 4).Value = "MarketMaking"
    $ws.Cells.Item(73, 5).Value = "DOWN"
    $ws.Cells.Item(73, 6).Value = 0.87
    $ws.Cells.Item(73, 7).Value = 0.9
    $ws.Cells.Item(73, 8).Value = "CLOSED"
    $ws.Cells.Item(73, 9).Value = 3.4483
    $ws.Cells.Item(73, 10).Value = 0.03
    $ws.Cells.Item(73, 11).Value = 100.56
    $ws.Cells.Item(73, 12).Value = 0
    $ws.Cells.Item(73, 13).Value = 0
    $ws.Cells.Item(73, 14).Value = 0.6
    $ws.Cells.Item(73, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item(73, 16).Value = "early_exit"
    $ws.Cells.Item(73, 17).Value = 0.13
}

Add-TradeRow73 $wb.Worksheets.Item("All Trades")
Add-TradeRow73 $wb.Worksheets.Item("MarketMaking")
